# Generate Report for Handoff
# Update "Latest Handoff Date/Datetime" values for the files that are
# being (re-)handed off: the one with "Handback transform failed" and
# the ones "Ready for handoff" (rows 4 and 6-10 on each sheet).

$wb = $excel.ActiveWorkbook

$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("D4").Value = "2016-38-17 03:38:10"
$ovw.Range("D6").Value = "2016-38-17 03:38:10"
$ovw.Range("D7").Value = "2016-38-17 03:38:10"
$ovw.Range("D8").Value = "2016-38-17 03:38:10"
$ovw.Range("D9").Value = "2016-38-17 03:38:10"
$ovw.Range("D10").Value = "2016-38-17 03:38:10"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-17 03:38:01"
$zhcn.Range("E6").Value = "2016-03-17 03:38:01"
$zhcn.Range("E7").Value = "2016-03-17 03:38:01"
$zhcn.Range("E8").Value = "2016-03-17 03:38:01"
$zhcn.Range("E9").Value = "2016-03-17 03:38:01"
$zhcn.Range("E10").Value = "2016-03-17 03:38:01"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-17 03:38:10"
$dede.Range("E6").Value = "2016-03-17 03:38:10"
$dede.Range("E7").Value = "2016-03-17 03:38:10"
$dede.Range("E8").Value = "2016-03-17 03:38:10"
$dede.Range("E9").Value = "2016-03-17 03:38:10"
$dede.Range("E10").Value = "2016-03-17 03:38:10"
